$d = $word.ActiveDocument

# --- Text replacements -----------------------------------------------

# Heading text: translate to Chinese
$d.Content.Find.Execute("Produce a table from regression results", $true, $false, $false, $false, $false,
                         $true, 1, $false, "线性回归结果", 2)

# Table header "fuel" -> "油耗"
$d.Content.Find.Execute("fuel", $true, $false, $false, $false, $false,
                         $true, 1, $false, "油耗", 2)

# Table row label "weight" -> "重量"
$d.Content.Find.Execute("weight", $true, $false, $false, $false, $false,
                         $true, 1, $false, "重量", 2)

# Numeric values in the weight row (coefficient, std. err., conf. interval)
$d.Content.Find.Execute(".001407", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".003102", 2)

$d.Content.Find.Execute(".0001008", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".0002223", 2)

$d.Content.Find.Execute(".001206", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".0026589", 2)

$d.Content.Find.Execute(".0016081", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".0035452", 2)

# --- Style cleanup -----------------------------------------------------
# Remove the custom "DocDefaults" paragraph style and detach "Normal" from it.

$normal = $d.Styles.Item("Normal")
$normal.BaseStyle = $null

$docDefaults = $d.Styles.Item("DocDefaults")
$docDefaults.Delete()
